# Updated cryptos list on Tue Nov 14 20:14:47 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.475.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.988.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.51"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.77"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.353"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0723"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.80%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.885"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.271.51"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.986.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.52%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.419.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.99"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -9.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.46%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.49%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.46"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.42"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.75%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -10.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0583"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0897"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -10.89%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -10.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.83"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.18"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.93%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0207"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.79%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0881"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.58%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.362.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.48"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.00%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.82%  "
